$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C2's value to "P@ssw0rd2"
$ws.Range("C2").Value = "P@ssw0rd2"

# Add hyperlink on C2 pointing to mailto:P@ssw0rd2 (mirrors existing C5 pattern)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:P@ssw0rd2")

# Update selection to D2
$ws.Range("D2").Select()
